$d = $word.ActiveDocument

# This edit removes the whole "Use Case" section that follows the project
# schedule table: a leading blank paragraph, the "Use Case" Heading2
# paragraph, the two-column "Use Case" details table, and the blank
# paragraph that trails it. The final blank paragraph (the one right
# before the section properties) is left untouched.

# Locate the two tables and the "Use Case" heading before making any
# modifications, since character offsets are most reliable when read
# prior to any write operation.
$table1 = $d.Tables.Item(1)
$useCaseTable = $d.Tables.Item(2)

$headingRange = $d.Content
$found = $headingRange.Find.Execute("Use Case", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Use Case' heading paragraph"
}
$headingStart = $headingRange.Start
$headingEnd = $headingRange.End

$leadingBlankStart = $table1.Range.End
$trailingBlankStart = $useCaseTable.Range.End

# Apply the deletions from the bottom of the document upward so that
# previously captured offsets for earlier content stay valid.

# 4. Blank paragraph left after the "Use Case" table.
$d.Range($trailingBlankStart, $trailingBlankStart + 1).Delete()

# 3. The "Use Case" details table itself.
$useCaseTable.Delete()

# 2. The "Use Case" Heading2 paragraph (its text plus the paragraph mark).
$d.Range($headingStart, $headingEnd + 1).Delete()

# 1. Blank paragraph that used to sit right after the schedule table.
$d.Range($leadingBlankStart, $leadingBlankStart + 1).Delete()
